$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 123.85714
$ws.Range("I53").Value = 59.7
$ws.Range("J53").Value = 284.25
$ws.Range("K53").Value = 59.7
$ws.Range("L53").Value = 284.25
$ws.Range("M53").Value = 577.3
$ws.Range("N53").Value = -1558.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 3729.9092
$ws.Range("I100").Value = 2218.7856
$ws.Range("K100").Value = 2218.7856
$ws.Range("M100").Value = -1677.7856

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 1366.5555
$ws.Range("I127").Value = 1299.8334
$ws.Range("J127").Value = 1500
$ws.Range("K127").Value = 3899.5002
$ws.Range("L127").Value = 4500
$ws.Range("M127").Value = 1060.4998
$ws.Range("N127").Value = -14420

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H130").Value = 47500
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 47500
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 47500
$ws.Range("M130").ClearContents()
$ws.Range("N130").Value = -57540

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5391.2173
$ws.Range("I138").Value = 2216.6667
$ws.Range("J138").Value = 6511.647
$ws.Range("K138").Value = 6650.000100000001
$ws.Range("L138").Value = 19534.941
$ws.Range("M138").Value = -1510.000100000001
$ws.Range("N138").Value = -29814.941

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 3290.1538
$ws.Range("I141").Value = 3023.5417
$ws.Range("K141").Value = 9070.625100000001
$ws.Range("M141").Value = -3890.625100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4967.081
$ws.Range("I32").Value = 4140.6875
$ws.Range("K32").Value = 4140.6875
$ws.Range("M32").Value = -3853.6875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 9990
$ws.Range("J19").Value = 9990
$ws.Range("L19").Value = 9990
$ws.Range("N19").Value = -10336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 60072.668
$ws.Range("J35").Value = 60072.668
$ws.Range("L35").Value = 60072.668
$ws.Range("N35").Value = -60692.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2567.9375
$ws.Range("I99").Value = 1275.9231
$ws.Range("K99").Value = 1275.9231
$ws.Range("M99").Value = 222.0769

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4995.4194
$ws.Range("I107").Value = 3515.913
$ws.Range("K107").Value = 3515.913
$ws.Range("M107").Value = -1595.913

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2590.35
$ws.Range("I122").Value = 984.3333
$ws.Range("J122").Value = 4999.375
$ws.Range("K122").Value = 2952.9999
$ws.Range("L122").Value = 14998.125
$ws.Range("M122").Value = -502.9998999999998
$ws.Range("N122").Value = -19898.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1562.6923
$ws.Range("I132").Value = 1562.6923
$ws.Range("K132").Value = 4688.0769
$ws.Range("M132").Value = -2158.0769

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2634
$ws.Range("I134").Value = 2767.7856
$ws.Range("J134").Value = 2366.4285
$ws.Range("K134").Value = 8303.356800000001
$ws.Range("L134").Value = 7099.2855
$ws.Range("M134").Value = -5768.356800000001
$ws.Range("N134").Value = -12169.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 80.375
$ws.Range("I2").Value = 53.875
$ws.Range("J2").Value = 106.875
$ws.Range("K2").Value = 323.25
$ws.Range("L2").Value = 641.25
$ws.Range("M2").Value = -210.25
$ws.Range("N2").Value = -867.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 106.117645
$ws.Range("J12").Value = 100.44444
$ws.Range("L12").Value = 301.33332
$ws.Range("N12").Value = -647.33332

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1185.5625
$ws.Range("I17").Value = 338.42856
$ws.Range("J17").Value = 1844.4445
$ws.Range("K17").Value = 1015.28568
$ws.Range("L17").Value = 5533.333500000001
$ws.Range("M17").Value = -846.28568
$ws.Range("N17").Value = -5871.333500000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 229.09091
$ws.Range("I33").Value = 245.71428
$ws.Range("J33").Value = 200
$ws.Range("K33").Value = 1474.28568
$ws.Range("L33").Value = 1200
$ws.Range("M33").Value = -1191.28568
$ws.Range("N33").Value = -1766

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H101").Value = 50000
$ws.Range("J101").Value = 50000
$ws.Range("L101").Value = 150000
$ws.Range("N101").Value = -154868

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 588
$ws.Range("I121").Value = 576.8
$ws.Range("K121").Value = 1730.4
$ws.Range("M121").Value = -420.3999999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1685.625
$ws.Range("J131").Value = 2018.25
$ws.Range("L131").Value = 6054.75
$ws.Range("N131").Value = -16134.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 6058.263
$ws.Range("I137").Value = 1932.2858
$ws.Range("J137").Value = 8465.083000000001
$ws.Range("K137").Value = 5796.857400000001
$ws.Range("L137").Value = 25395.249
$ws.Range("M137").Value = -696.8574000000008
$ws.Range("N137").Value = -35595.249

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1652.75
$ws.Range("I122").Value = 1537
$ws.Range("K122").Value = 4611
$ws.Range("M122").Value = -2161

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 7281
$ws.Range("I126").Value = 3117
$ws.Range("J126").Value = 12833
$ws.Range("K126").Value = 9351
$ws.Range("L126").Value = 38499
$ws.Range("M126").Value = -6881
$ws.Range("N126").Value = -43439

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2240.125
$ws.Range("I132").Value = 2320.3333
$ws.Range("J132").Value = 1999.5
$ws.Range("K132").Value = 6960.999899999999
$ws.Range("L132").Value = 5998.5
$ws.Range("M132").Value = -4430.999899999999
$ws.Range("N132").Value = -11058.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 75000
$ws.Range("J135").Value = 75000
$ws.Range("L135").Value = 75000
$ws.Range("N135").Value = -85140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2512.6428
$ws.Range("I132").Value = 1969.7727
$ws.Range("J132").Value = 4503.1665
$ws.Range("K132").Value = 5909.3181
$ws.Range("L132").Value = 13509.4995
$ws.Range("M132").Value = -3379.3181
$ws.Range("N132").Value = -18569.4995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H139").Value = 109107.5
$ws.Range("J139").Value = 109107.5
$ws.Range("L139").Value = 109107.5
$ws.Range("N139").Value = -119387.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 150250
$ws.Range("I19").Value = 150250
$ws.Range("K19").Value = 150250
$ws.Range("M19").Value = -150076

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1576.5555
$ws.Range("I81").Value = 1576.5555
$ws.Range("K81").Value = 3153.111
$ws.Range("M81").Value = -2092.111

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1576.5555
$ws.Range("I84").Value = 1576.5555
$ws.Range("K84").Value = 15765.555
$ws.Range("M84").Value = -10461.555

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1473.909
$ws.Range("J113").Value = 1357
$ws.Range("L113").Value = 4071
$ws.Range("N113").Value = -8411

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4080.2415
$ws.Range("I132").Value = 2995.3
$ws.Range("J132").Value = 4651.263
$ws.Range("K132").Value = 8985.900000000001
$ws.Range("L132").Value = 13953.789
$ws.Range("M132").Value = -6455.900000000001
$ws.Range("N132").Value = -19013.789
